$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (status) for the "need to clean up piano keys a little" row (row 8)
# changes from "Pending" to "Resolved" since the columns no longer depend on
# the starting key.
$ws.Range("C8").Value = "Resolved"

# Add a note in column D for row 8.
$ws.Range("D8").Value = "draft58"

# Update the active selection to reflect the last edited cell.
$ws.Range("D8").Select()
